$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows dropped from the 2024 Sudan (coa) series: Djibouti, Turkiye, Uganda.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(506).Delete()  # Uganda
$ws.Rows.Item(505).Delete()  # Turkiye
$ws.Rows.Item(492).Delete()  # Djibouti

# Refresh the short-url token used across every data row (column B).
$ws.Range("B2:B504").Value = "iqQSa0"

# Update the refreshed 2024 statistics (columns D, N, O, P, Q, T, V) for the remaining coo rows.
# 486: Burundi
$ws.Range("D486").Value = "485"
$ws.Range("N486").Value = "8"
$ws.Range("O486").Value = "7"
$ws.Range("P486").Value = "0"
$ws.Range("Q486").Value = "0"
$ws.Range("T486").Value = "0"
$ws.Range("V486").Value = "0"
# 487: Central African Rep.
$ws.Range("D487").Value = "486"
$ws.Range("N487").Value = "10053"
$ws.Range("O487").Value = "5"
$ws.Range("P487").Value = "1565"
$ws.Range("Q487").Value = "0"
$ws.Range("T487").Value = "0"
$ws.Range("V487").Value = "0"
# 488: Chad
$ws.Range("D488").Value = "487"
$ws.Range("N488").Value = "955"
$ws.Range("O488").Value = "6"
$ws.Range("P488").Value = "0"
$ws.Range("Q488").Value = "0"
$ws.Range("T488").Value = "7"
$ws.Range("V488").Value = "0"
# 489: Cameroon
$ws.Range("D489").Value = "488"
$ws.Range("N489").Value = "7"
$ws.Range("O489").Value = "0"
$ws.Range("P489").Value = "0"
$ws.Range("Q489").Value = "0"
$ws.Range("T489").Value = "0"
$ws.Range("V489").Value = "0"
# 490: Congo
$ws.Range("D490").Value = "489"
$ws.Range("N490").Value = "10"
$ws.Range("O490").Value = "5"
$ws.Range("P490").Value = "0"
$ws.Range("Q490").Value = "0"
$ws.Range("T490").Value = "0"
$ws.Range("V490").Value = "0"
# 491: Dem. Rep. of the Congo
$ws.Range("D491").Value = "490"
$ws.Range("N491").Value = "331"
$ws.Range("O491").Value = "260"
$ws.Range("P491").Value = "0"
$ws.Range("Q491").Value = "0"
$ws.Range("T491").Value = "0"
$ws.Range("V491").Value = "0"
# 492: Eritrea
$ws.Range("D492").Value = "491"
$ws.Range("N492").Value = "108349"
$ws.Range("O492").Value = "23664"
$ws.Range("P492").Value = "0"
$ws.Range("Q492").Value = "0"
$ws.Range("T492").Value = "49"
$ws.Range("V492").Value = "0"
# 493: Ethiopia
$ws.Range("D493").Value = "492"
$ws.Range("N493").Value = "49547"
$ws.Range("O493").Value = "20348"
$ws.Range("P493").Value = "11663"
$ws.Range("Q493").Value = "0"
$ws.Range("T493").Value = "3289"
$ws.Range("V493").Value = "0"
# 494: Palestinian
$ws.Range("D494").Value = "493"
$ws.Range("N494").Value = "49"
$ws.Range("O494").Value = "24"
$ws.Range("P494").Value = "0"
$ws.Range("Q494").Value = "0"
$ws.Range("T494").Value = "187"
$ws.Range("V494").Value = "0"
# 495: Iraq
$ws.Range("D495").Value = "494"
$ws.Range("N495").Value = "9"
$ws.Range("O495").Value = "13"
$ws.Range("P495").Value = "0"
$ws.Range("Q495").Value = "0"
$ws.Range("T495").Value = "5"
$ws.Range("V495").Value = "0"
# 496: Jordan
$ws.Range("D496").Value = "495"
$ws.Range("N496").Value = "0"
$ws.Range("O496").Value = "0"
$ws.Range("P496").Value = "0"
$ws.Range("Q496").Value = "0"
$ws.Range("T496").Value = "6"
$ws.Range("V496").Value = "0"
# 497: Kuwait
$ws.Range("D497").Value = "496"
$ws.Range("N497").Value = "0"
$ws.Range("O497").Value = "0"
$ws.Range("P497").Value = "0"
$ws.Range("Q497").Value = "0"
$ws.Range("T497").Value = "5"
$ws.Range("V497").Value = "0"
# 498: Nigeria
$ws.Range("D498").Value = "497"
$ws.Range("N498").Value = "12"
$ws.Range("O498").Value = "0"
$ws.Range("P498").Value = "387"
$ws.Range("Q498").Value = "0"
$ws.Range("T498").Value = "0"
$ws.Range("V498").Value = "0"
# 499: Saudi Arabia
$ws.Range("D499").Value = "498"
$ws.Range("N499").Value = "0"
$ws.Range("O499").Value = "0"
$ws.Range("P499").Value = "0"
$ws.Range("Q499").Value = "0"
$ws.Range("T499").Value = "5"
$ws.Range("V499").Value = "0"
# 500: Somalia
$ws.Range("D500").Value = "499"
$ws.Range("N500").Value = "258"
$ws.Range("O500").Value = "77"
$ws.Range("P500").Value = "33"
$ws.Range("Q500").Value = "0"
$ws.Range("T500").Value = "0"
$ws.Range("V500").Value = "0"
# 501: South Sudan
$ws.Range("D501").Value = "500"
$ws.Range("N501").Value = "613052"
$ws.Range("O501").Value = "0"
$ws.Range("P501").Value = "292919"
$ws.Range("Q501").Value = "0"
$ws.Range("T501").Value = "0"
$ws.Range("V501").Value = "0"
# 502: Sudan
$ws.Range("D502").Value = "501"
$ws.Range("N502").Value = "0"
$ws.Range("O502").Value = "0"
$ws.Range("P502").Value = "0"
$ws.Range("Q502").Value = "11559970"
$ws.Range("T502").Value = "0"
$ws.Range("V502").Value = "32823"
# 503: Syrian Arab Rep.
$ws.Range("D503").Value = "502"
$ws.Range("N503").Value = "9860"
$ws.Range("O503").Value = "0"
$ws.Range("P503").Value = "0"
$ws.Range("Q503").Value = "0"
$ws.Range("T503").Value = "0"
$ws.Range("V503").Value = "0"
# 504: Yemen
$ws.Range("D504").Value = "503"
$ws.Range("N504").Value = "877"
$ws.Range("O504").Value = "173"
$ws.Range("P504").Value = "0"
$ws.Range("Q504").Value = "0"
$ws.Range("T504").Value = "0"
$ws.Range("V504").Value = "0"
